$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Autonomous_temporary")

# --- Insert two new rows before the "Message: ACU_status" block (old row 82) ---
# This shifts the old rows 82-86 down to 84-88, while row 81 (previously blank)
# stays in place and gets new content below.
$ws.Rows("82:83").Insert()

# Stamp the "data row" style (border, normal font/fill - style index 3) onto the
# two freshly inserted rows by copying the format from an existing data row.
$ws.Range("A80:K80").Copy()
$ws.Range("A81:K83").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 81: R2D_button_raw signal
$ws.Range("A81").Value = "R2D_button_raw"
$ws.Range("B81").Value = 48
$ws.Range("C81").Value = 8
$ws.Range("D81").Value = "Intel"
$ws.Range("E81").Value = $false
$ws.Range("F81").Value = 1
$ws.Range("G81").Value = 0

# Row 82: Ignition_switch_raw signal
$ws.Range("A82").Value = "Ignition_switch_raw"
$ws.Range("B82").Value = 56
$ws.Range("C82").Value = 8
$ws.Range("D82").Value = "Intel"
$ws.Range("E82").Value = $false
$ws.Range("F82").Value = 1
$ws.Range("G82").Value = 0

# Row 83 stays blank (just an inserted spacer row, matches old row 81's spacer role).

# --- Append the new "Message: VCU_APPS_RAW" block after the ACU_status block ---
# Old sheet ended at row 86 (Internal_temperature); after the 2-row insert above
# that content now sits at row 88, so the new block starts two rows later at 90.

# Message header row (style 1: bold + blue fill) - copy format from the untouched
# row 1 header, which carries the same style.
$ws.Range("A1:C1").Copy()
$ws.Range("A90:C90").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A90").Value = "Message: VCU_APPS_RAW"
$ws.Range("B90").Value = "ID: 0x610"
$ws.Range("C90").Value = "Sender(s): VCU"

# Column header row (style 2: bold + yellow fill) - copy format from the untouched
# row 2 header, which carries the same style.
$ws.Range("A2:K2").Copy()
$ws.Range("A91:K91").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("A91").Value = "Signal Name"
$ws.Range("B91").Value = "Start Bit"
$ws.Range("C91").Value = "Length (bits)"
$ws.Range("D91").Value = "Byte Order"
$ws.Range("E91").Value = "Signed"
$ws.Range("F91").Value = "Factor"
$ws.Range("G91").Value = "Offset"
$ws.Range("H91").Value = "Min"
$ws.Range("I91").Value = "Max"
$ws.Range("J91").Value = "Unit"
$ws.Range("K91").Value = "Choices"

# Data rows (style 3: bordered data row)
$ws.Range("A81:K81").Copy()
$ws.Range("A92:K95").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

# Row 92: APPS_1_raw_bits
$ws.Range("A92").Value = "APPS_1_raw_bits"
$ws.Range("B92").Value = 0
$ws.Range("C92").Value = 8
$ws.Range("D92").Value = "Intel"
$ws.Range("E92").Value = $false
$ws.Range("F92").Value = 1
$ws.Range("G92").Value = 0

# Row 93: APPS_2_raw_bits
$ws.Range("A93").Value = "APPS_2_raw_bits"
$ws.Range("B93").Value = 8
$ws.Range("C93").Value = 8
$ws.Range("D93").Value = "Intel"
$ws.Range("E93").Value = $false
$ws.Range("F93").Value = 1
$ws.Range("G93").Value = 0

# Row 94: delta_raw
$ws.Range("A94").Value = "delta_raw"
$ws.Range("B94").Value = 16
$ws.Range("C94").Value = 8
$ws.Range("D94").Value = "Intel"
$ws.Range("E94").Value = $false
$ws.Range("F94").Value = 1
$ws.Range("G94").Value = 0

# Row 95: cpu_temp
$ws.Range("A95").Value = "cpu_temp"
$ws.Range("B95").Value = 24
$ws.Range("C95").Value = 8
$ws.Range("D95").Value = "Intel"
$ws.Range("E95").Value = $false
$ws.Range("F95").Value = 1
$ws.Range("G95").Value = 0
